$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first AV1 advance-related line (was 002/RRR) ---
$ws.Range("A2").Value = "001/RRR/AV1"
$ws.Range("C2").Value = "B219321"
$ws.Range("D2").Value = "JEMAA HORMI"
$ws.Range("H2").Value = "--"
$ws.Range("J2").Value = "--"
$ws.Range("L2").Value = 7000
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 7000

# --- Row 3: second AV1 advance-related line (was 001/TTT) ---
$ws.Range("A3").Value = "001/RRR/AV1"
$ws.Range("C3").Value = "I83603"
$ws.Range("D3").Value = "MOHAMED BADRANE"
$ws.Range("H3").Value = "--"
$ws.Range("J3").Value = "--"
$ws.Range("L3").Value = 14000
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 14000

# --- Insert two new data rows before the old total row (old row 4 -> new row 6) ---
$ws.Rows("4:5").Insert()

# --- Row 4: new data row for JEMAA HORMI ---
$ws.Range("A4").Value = "001/RRR/AV1"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "B219321"
$ws.Range("D4").Value = "JEMAA HORMI"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 1000

# --- Row 5: new data row for MOHAMED BADRANE ---
$ws.Range("A5").Value = "001/RRR/AV1"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "I83603"
$ws.Range("D5").Value = "MOHAMED BADRANE"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 2000

# --- Row 6: totals row (previously row 4, shifted down by inserting 2 rows) ---
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 21000
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 24000
